# Update the 20x5 answer table in place: each cell keeps its run/paragraph
# formatting (rFonts TimeNewRoman, sz 30) intact because we only rewrite the
# cell Range.Text, not the cell/paragraph/run structure itself.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "5+35=40"
$t.Cell(1,2).Range.Text = "83-36=47"
$t.Cell(1,3).Range.Text = "1+48=49"
$t.Cell(1,4).Range.Text = "43+48=91"
$t.Cell(1,5).Range.Text = "38-15=23"

$t.Cell(2,1).Range.Text = "93-67=26"
$t.Cell(2,2).Range.Text = "15-0=15"
$t.Cell(2,3).Range.Text = "94+4=98"
$t.Cell(2,4).Range.Text = "42-2=40"
$t.Cell(2,5).Range.Text = "91-18=73"

$t.Cell(3,1).Range.Text = "97-50=47"
$t.Cell(3,2).Range.Text = "79-29=50"
$t.Cell(3,3).Range.Text = "39+58=97"
$t.Cell(3,4).Range.Text = "63-12=51"
$t.Cell(3,5).Range.Text = "28-19=9"

$t.Cell(4,1).Range.Text = "49-34=15"
$t.Cell(4,2).Range.Text = "1+96=97"
$t.Cell(4,3).Range.Text = "91-27=64"
$t.Cell(4,4).Range.Text = "56-17=39"
$t.Cell(4,5).Range.Text = "77-67=10"

$t.Cell(5,1).Range.Text = "40-28=12"
$t.Cell(5,2).Range.Text = "94-57=37"
$t.Cell(5,3).Range.Text = "68+14=82"
$t.Cell(5,4).Range.Text = "8+84=92"
$t.Cell(5,5).Range.Text = "65+29=94"

$t.Cell(6,1).Range.Text = "49-16=33"
$t.Cell(6,2).Range.Text = "6+39=45"
$t.Cell(6,3).Range.Text = "84-23=61"
$t.Cell(6,4).Range.Text = "92+6=98"
$t.Cell(6,5).Range.Text = "21+41=62"

$t.Cell(7,1).Range.Text = "48+16=64"
$t.Cell(7,2).Range.Text = "25+25=50"
$t.Cell(7,3).Range.Text = "62-40=22"
$t.Cell(7,4).Range.Text = "29+11=40"
$t.Cell(7,5).Range.Text = "28+20=48"

$t.Cell(8,1).Range.Text = "72-0=72"
$t.Cell(8,2).Range.Text = "88-21=67"
$t.Cell(8,3).Range.Text = "84+3=87"
$t.Cell(8,4).Range.Text = "97-71=26"
$t.Cell(8,5).Range.Text = "61-19=42"

$t.Cell(9,1).Range.Text = "4+83=87"
$t.Cell(9,2).Range.Text = "7+2=9"
$t.Cell(9,3).Range.Text = "54+27=81"
$t.Cell(9,4).Range.Text = "76+21=97"
$t.Cell(9,5).Range.Text = "77-60=17"

$t.Cell(10,1).Range.Text = "11+18=29"
$t.Cell(10,2).Range.Text = "40+32=72"
$t.Cell(10,3).Range.Text = "60+3=63"
$t.Cell(10,4).Range.Text = "62-23=39"
$t.Cell(10,5).Range.Text = "47+0=47"

$t.Cell(11,1).Range.Text = "56-25=31"
$t.Cell(11,2).Range.Text = "99-64=35"
$t.Cell(11,3).Range.Text = "57-4=53"
$t.Cell(11,4).Range.Text = "3+2=5"
$t.Cell(11,5).Range.Text = "29+61=90"

$t.Cell(12,1).Range.Text = "61-44=17"
$t.Cell(12,2).Range.Text = "60-4=56"
$t.Cell(12,3).Range.Text = "43-28=15"
$t.Cell(12,4).Range.Text = "50+35=85"
$t.Cell(12,5).Range.Text = "92-46=46"

$t.Cell(13,1).Range.Text = "40-38=2"
$t.Cell(13,2).Range.Text = "82-30=52"
$t.Cell(13,3).Range.Text = "5+91=96"
$t.Cell(13,4).Range.Text = "49+34=83"
$t.Cell(13,5).Range.Text = "45+33=78"

$t.Cell(14,1).Range.Text = "46+14=60"
$t.Cell(14,2).Range.Text = "68+18=86"
$t.Cell(14,3).Range.Text = "35-17=18"
$t.Cell(14,4).Range.Text = "3+47=50"
$t.Cell(14,5).Range.Text = "43-12=31"

$t.Cell(15,1).Range.Text = "74-65=9"
$t.Cell(15,2).Range.Text = "5+17=22"
$t.Cell(15,3).Range.Text = "0+38=38"
$t.Cell(15,4).Range.Text = "59-5=54"
$t.Cell(15,5).Range.Text = "41+56=97"

$t.Cell(16,1).Range.Text = "79-34=45"
$t.Cell(16,2).Range.Text = "76-38=38"
$t.Cell(16,3).Range.Text = "77-4=73"
$t.Cell(16,4).Range.Text = "97-24=73"
$t.Cell(16,5).Range.Text = "85-27=58"

$t.Cell(17,1).Range.Text = "50+41=91"
$t.Cell(17,2).Range.Text = "54+42=96"
$t.Cell(17,3).Range.Text = "74-28=46"
$t.Cell(17,4).Range.Text = "30-22=8"
$t.Cell(17,5).Range.Text = "2+93=95"

$t.Cell(18,1).Range.Text = "90-25=65"
$t.Cell(18,2).Range.Text = "87+7=94"
$t.Cell(18,3).Range.Text = "62+21=83"
$t.Cell(18,4).Range.Text = "31+10=41"
$t.Cell(18,5).Range.Text = "90-25=65"

$t.Cell(19,1).Range.Text = "5+22=27"
$t.Cell(19,2).Range.Text = "0+57=57"
$t.Cell(19,3).Range.Text = "13+79=92"
$t.Cell(19,4).Range.Text = "27-11=16"
$t.Cell(19,5).Range.Text = "46-5=41"

$t.Cell(20,1).Range.Text = "64-31=33"
$t.Cell(20,2).Range.Text = "53-10=43"
$t.Cell(20,3).Range.Text = "34+37=71"
$t.Cell(20,4).Range.Text = "22+12=34"
$t.Cell(20,5).Range.Text = "69+11=80"
